$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update probability matrix values per "games pulled march 7" refresh.
# Diff only touches value cells (row 2 through row 19) on the single sheet.
$ws.Cells.Item(2, 2).Value = 0.1862244897959184
$ws.Cells.Item(2, 3).Value = 0.576530612244898
$ws.Cells.Item(2, 10).Value = 0.01275510204081633
$ws.Cells.Item(2, 15).Value = 0.002551020408163265
$ws.Cells.Item(2, 16).Value = 0.1428571428571428
$ws.Cells.Item(2, 19).Value = 0.07908163265306123
$ws.Cells.Item(3, 2).Value = 0.008583690987124463
$ws.Cells.Item(3, 3).Value = 0.02575107296137339
$ws.Cells.Item(3, 10).Value = 0.03004291845493562
$ws.Cells.Item(3, 16).Value = 0.7725321888412017
$ws.Cells.Item(3, 19).Value = 0.1630901287553648
$ws.Cells.Item(4, 10).Value = 0.02
$ws.Cells.Item(4, 16).Value = 0.68
$ws.Cells.Item(4, 19).Value = 0.3
$ws.Cells.Item(6, 2).Value = 0.06666666666666667
$ws.Cells.Item(6, 4).Value = 0.01333333333333333
$ws.Cells.Item(6, 6).Value = 0.12
$ws.Cells.Item(6, 10).Value = 0.1777777777777778
$ws.Cells.Item(6, 15).Value = 0.04888888888888889
$ws.Cells.Item(6, 17).Value = 0.1511111111111111
$ws.Cells.Item(6, 18).Value = 0.06666666666666667
$ws.Cells.Item(6, 19).Value = 0.3555555555555556
$ws.Cells.Item(7, 2).Value = 0.1471698113207547
$ws.Cells.Item(7, 4).Value = 0.02264150943396226
$ws.Cells.Item(7, 6).Value = 0.0339622641509434
$ws.Cells.Item(7, 10).Value = 0.1245283018867925
$ws.Cells.Item(7, 15).Value = 0.007547169811320755
$ws.Cells.Item(7, 17).Value = 0.1283018867924528
$ws.Cells.Item(7, 18).Value = 0.06415094339622641
$ws.Cells.Item(7, 19).Value = 0.4716981132075472
$ws.Cells.Item(8, 2).Value = 0.09652509652509653
$ws.Cells.Item(8, 4).Value = 0.01737451737451737
$ws.Cells.Item(8, 6).Value = 0.06177606177606178
$ws.Cells.Item(8, 10).Value = 0.1177606177606178
$ws.Cells.Item(8, 15).Value = 0.01737451737451737
$ws.Cells.Item(8, 18).Value = 0.07722007722007722
$ws.Cells.Item(8, 19).Value = 0.4691119691119691
$ws.Cells.Item(9, 2).Value = 0.1181818181818182
$ws.Cells.Item(9, 4).Value = 0.01363636363636364
$ws.Cells.Item(9, 6).Value = 0.08636363636363636
$ws.Cells.Item(9, 10).Value = 0.1636363636363636
$ws.Cells.Item(9, 15).Value = 0.004545454545454545
$ws.Cells.Item(9, 17).Value = 0.1454545454545454
$ws.Cells.Item(9, 18).Value = 0.08181818181818182
$ws.Cells.Item(9, 19).Value = 0.3863636363636364
$ws.Cells.Item(10, 2).Value = 0.1420807453416149
$ws.Cells.Item(10, 4).Value = 0.02329192546583851
$ws.Cells.Item(10, 6).Value = 0.06055900621118013
$ws.Cells.Item(10, 10).Value = 0.1583850931677019
$ws.Cells.Item(10, 15).Value = 0.01940993788819876
$ws.Cells.Item(10, 17).Value = 0.1801242236024845
$ws.Cells.Item(10, 18).Value = 0.06521739130434782
$ws.Cells.Item(10, 19).Value = 0.3509316770186335
$ws.Cells.Item(11, 7).Value = 0.1295774647887324
$ws.Cells.Item(11, 10).Value = 0.07323943661971831
$ws.Cells.Item(11, 11).Value = 0.171830985915493
$ws.Cells.Item(11, 12).Value = 0.6084507042253521
$ws.Cells.Item(11, 19).Value = 0.01690140845070422
$ws.Cells.Item(12, 7).Value = 0.7835497835497836
$ws.Cells.Item(12, 10).Value = 0.1471861471861472
$ws.Cells.Item(12, 11).Value = 0.004329004329004329
$ws.Cells.Item(12, 12).Value = 0.03896103896103896
$ws.Cells.Item(12, 19).Value = 0.02597402597402598
$ws.Cells.Item(13, 7).Value = 0.6857142857142857
$ws.Cells.Item(13, 10).Value = 0.2571428571428571
$ws.Cells.Item(13, 19).Value = 0.05714285714285714
$ws.Cells.Item(14, 7).Value = 0.4
$ws.Cells.Item(14, 10).Value = 0.2
$ws.Cells.Item(14, 19).Value = 0.4
$ws.Cells.Item(15, 6).Value = 0.02521008403361345
$ws.Cells.Item(15, 8).Value = 0.1974789915966386
$ws.Cells.Item(15, 9).Value = 0.07142857142857142
$ws.Cells.Item(15, 10).Value = 0.2394957983193277
$ws.Cells.Item(15, 11).Value = 0.1008403361344538
$ws.Cells.Item(15, 13).Value = 0.02941176470588235
$ws.Cells.Item(15, 14).Value = 0.004201680672268907
$ws.Cells.Item(15, 15).Value = 0.09663865546218488
$ws.Cells.Item(15, 19).Value = 0.2352941176470588
$ws.Cells.Item(16, 6).Value = 0.02334630350194553
$ws.Cells.Item(16, 8).Value = 0.2373540856031128
$ws.Cells.Item(16, 9).Value = 0.05836575875486381
$ws.Cells.Item(16, 10).Value = 0.3073929961089494
$ws.Cells.Item(16, 11).Value = 0.132295719844358
$ws.Cells.Item(16, 13).Value = 0.0311284046692607
$ws.Cells.Item(16, 15).Value = 0.06614785992217899
$ws.Cells.Item(16, 19).Value = 0.1439688715953307
$ws.Cells.Item(17, 6).Value = 0.01237623762376238
$ws.Cells.Item(17, 8).Value = 0.1633663366336634
$ws.Cells.Item(17, 9).Value = 0.1014851485148515
$ws.Cells.Item(17, 10).Value = 0.4381188118811881
$ws.Cells.Item(17, 11).Value = 0.08663366336633663
$ws.Cells.Item(17, 13).Value = 0.03217821782178218
$ws.Cells.Item(17, 14).Value = 0.002475247524752475
$ws.Cells.Item(17, 15).Value = 0.0594059405940594
$ws.Cells.Item(17, 19).Value = 0.103960396039604
$ws.Cells.Item(18, 6).Value = 0.01694915254237288
$ws.Cells.Item(18, 8).Value = 0.1751412429378531
$ws.Cells.Item(18, 9).Value = 0.096045197740113
$ws.Cells.Item(18, 10).Value = 0.4011299435028249
$ws.Cells.Item(18, 11).Value = 0.0903954802259887
$ws.Cells.Item(18, 13).Value = 0.01129943502824859
$ws.Cells.Item(18, 14).Value = 0.005649717514124294
$ws.Cells.Item(18, 15).Value = 0.06779661016949153
$ws.Cells.Item(18, 19).Value = 0.1355932203389831
$ws.Cells.Item(19, 6).Value = 0.0134847409510291
$ws.Cells.Item(19, 8).Value = 0.2178850248403123
$ws.Cells.Item(19, 9).Value = 0.09013484740951029
$ws.Cells.Item(19, 10).Value = 0.3307310149041874
$ws.Cells.Item(19, 11).Value = 0.1298793470546487
$ws.Cells.Item(19, 13).Value = 0.028388928317956
$ws.Cells.Item(19, 14).Value = 0.0021291696238467
$ws.Cells.Item(19, 15).Value = 0.05748757984386089
$ws.Cells.Item(19, 19).Value = 0.1298793470546487

$wb.Save()
